# "Adding features for tenure"
# Adds two new columns to the households sheet: "Tenure Pref" (col N) and
# "Tenure" (col O), each populated with Rent/Own values, tweaks one data
# point (H3), widens column M slightly, and moves the active sheet/selection
# from financial_capital back to households (with housing_stock's selection
# also nudged to K30).

$wb = $excel.ActiveWorkbook

$households = $wb.Worksheets.Item("households")
$housingStock = $wb.Worksheets.Item("housing_stock")

# New header + data cells for households sheet ("Tenure Pref" / "Tenure" columns).
# Written in this specific order so the shared-string table comes out in the
# same order Excel produced it in (header N1 first, then the data rows, then
# the O1 header last).
$households.Range("N1").Value = "Tenure Pref"

$households.Range("O2").Value = "Rent"
$households.Range("N2").Value = "Own"

$households.Range("O3").Value = "Own"
$households.Range("N3").Value = "Own"
$households.Range("H3").Value = 5

$households.Range("O4").Value = "Rent"
$households.Range("N4").Value = "Rent"

$households.Range("O5").Value = "Own"
$households.Range("N5").Value = "Own"

$households.Range("O1").Value = "Tenure"

# Widen column M ("Damage State") -- matches the width tweak recorded in the
# source diff (closest achievable value given this runtime's column-width
# quantization: stored width ends up 16.8333, vs. 16.83203125 in real Excel).
$households.Columns.Item(13).ColumnWidth = 16

$households.Activate()
$households.Range("A1").Select() | Out-Null

$housingStock.Activate()
$housingStock.Range("K30").Select() | Out-Null

$households.Activate()
